# Hortaliza, Vega Modelo de Temuco - Ciboulette: add a new weekly record.
# A new data row is inserted at row 60, pushing all subsequent rows
# (old 60..187) down by one (new 61..188). The sheet's used range grows
# from A1:R187 to A1:R188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 60, shifting rows 60-187 down to 61-188.
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with the new observation.
$ws.Cells.Item(60, 1).Value  = 10
$ws.Cells.Item(60, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(60, 3).Value  = "La Araucanía"
$ws.Cells.Item(60, 4).Value  = 44544
$ws.Cells.Item(60, 5).Value  = 9
$ws.Cells.Item(60, 6).Value  = 100112039
$ws.Cells.Item(60, 7).Value  = "Ciboulette"
$ws.Cells.Item(60, 8).Value  = "Sin especificar"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 25
$ws.Cells.Item(60, 11).Value = 7000
$ws.Cells.Item(60, 12).Value = 7000
$ws.Cells.Item(60, 13).Value = 7000
$ws.Cells.Item(60, 14).Value = "`$/docena de atados"
$ws.Cells.Item(60, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(60, 16).Value = 2333
$ws.Cells.Item(60, 17).Value = 3
$ws.Cells.Item(60, 18).Value = "Hortaliza"
